$wb = $excel.ActiveWorkbook

# ALC row 33
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(33, 8).Value = 129.03847
$ws.Cells.Item(33, 9).Value = 81.85714
$ws.Cells.Item(33, 10).Value = 327.2
$ws.Cells.Item(33, 11).Value = 81.85714
$ws.Cells.Item(33, 12).Value = 327.2
$ws.Cells.Item(33, 13).Value = 147.14286
$ws.Cells.Item(33, 14).Value = -785.2

# ALC row 137
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(137, 8).Value = 1285.5555
$ws.Cells.Item(137, 9).Value = 788.375
$ws.Cells.Item(137, 10).Value = 1683.3
$ws.Cells.Item(137, 11).Value = 2365.125
$ws.Cells.Item(137, 12).Value = 5049.9
$ws.Cells.Item(137, 13).Value = 184.875
$ws.Cells.Item(137, 14).Value = -10149.9

# ALC row 138
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(138, 8).Value = 9436405
$ws.Cells.Item(138, 9).Value = 2426.04
$ws.Cells.Item(138, 10).Value = 17859600
$ws.Cells.Item(138, 11).Value = 7278.12
$ws.Cells.Item(138, 12).Value = 53578800
$ws.Cells.Item(138, 13).Value = -2138.12
$ws.Cells.Item(138, 14).Value = -53589080

# ARM row 74
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(74, 8).Value = 976.5625
$ws.Cells.Item(74, 9).Value = 931.25
$ws.Cells.Item(74, 11).Value = 931.25
$ws.Cells.Item(74, 13).Value = -57.25

# ARM row 77
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(77, 8).Value = 976.5625
$ws.Cells.Item(77, 9).Value = 931.25
$ws.Cells.Item(77, 11).Value = 4656.25
$ws.Cells.Item(77, 13).Value = -288.25

# BSM row 86
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(86, 8).Value = 23758.143
$ws.Cells.Item(86, 9).Value = 15075
$ws.Cells.Item(86, 10).Value = 35335.668
$ws.Cells.Item(86, 11).Value = 15075
$ws.Cells.Item(86, 12).Value = 35335.668
$ws.Cells.Item(86, 13).Value = -13952
$ws.Cells.Item(86, 14).Value = -37581.668

# BSM row 89
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(89, 8).Value = 23758.143
$ws.Cells.Item(89, 9).Value = 15075
$ws.Cells.Item(89, 10).Value = 35335.668
$ws.Cells.Item(89, 11).Value = 75375
$ws.Cells.Item(89, 12).Value = 176678.34
$ws.Cells.Item(89, 13).Value = -69759
$ws.Cells.Item(89, 14).Value = -187910.34

# BSM row 134
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(134, 8).Value = 2244.9592
$ws.Cells.Item(134, 9).Value = 1297.7297
$ws.Cells.Item(134, 10).Value = 5165.5835
$ws.Cells.Item(134, 11).Value = 3893.189100000001
$ws.Cells.Item(134, 12).Value = 15496.7505
$ws.Cells.Item(134, 13).Value = -1358.189100000001
$ws.Cells.Item(134, 14).Value = -20566.7505

# CRP row 58
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(58, 8).Value = 1261.7333
$ws.Cells.Item(58, 9).Value = 381.1
$ws.Cells.Item(58, 10).Value = 3023
$ws.Cells.Item(58, 11).Value = 381.1
$ws.Cells.Item(58, 12).Value = 3023
$ws.Cells.Item(58, 13).Value = -178.1
$ws.Cells.Item(58, 14).Value = -3429

# CRP row 87
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(87, 8).Value = 800330
$ws.Cells.Item(87, 10).Value = 800330
$ws.Cells.Item(87, 12).Value = 800330
$ws.Cells.Item(87, 14).Value = -802702

# CRP row 90
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(90, 8).Value = 800330
$ws.Cells.Item(90, 10).Value = 800330
$ws.Cells.Item(90, 12).Value = 2400990
$ws.Cells.Item(90, 14).Value = -2412846

# CRP row 94
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(94, 8).Value = 903.9394
$ws.Cells.Item(94, 9).Value = 903.875
$ws.Cells.Item(94, 10).Value = 903.96
$ws.Cells.Item(94, 11).Value = 903.875
$ws.Cells.Item(94, 12).Value = 903.96
$ws.Cells.Item(94, 13).Value = -452.875
$ws.Cells.Item(94, 14).Value = -1805.96

# CRP row 118
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(118, 8).Value = 30000
$ws.Cells.Item(118, 10).Value = 30000
$ws.Cells.Item(118, 12).Value = 30000
$ws.Cells.Item(118, 14).Value = -33314

# CRP row 136
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(136, 8).Value = 1261.7333
$ws.Cells.Item(136, 9).Value = 381.1
$ws.Cells.Item(136, 10).Value = 3023
$ws.Cells.Item(136, 11).Value = 1143.3
$ws.Cells.Item(136, 12).Value = 9069
$ws.Cells.Item(136, 13).Value = 1406.7
$ws.Cells.Item(136, 14).Value = -14169

# CUL row 68
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(68, 8).Value = 849.11
$ws.Cells.Item(68, 9).Value = 726.95386
$ws.Cells.Item(68, 10).Value = 1075.9714
$ws.Cells.Item(68, 11).Value = 2180.86158
$ws.Cells.Item(68, 12).Value = 3227.9142
$ws.Cells.Item(68, 13).Value = -1369.86158
$ws.Cells.Item(68, 14).Value = -4849.914199999999

# CUL row 71
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(71, 8).Value = 849.11
$ws.Cells.Item(71, 9).Value = 726.95386
$ws.Cells.Item(71, 10).Value = 1075.9714
$ws.Cells.Item(71, 11).Value = 6542.584739999999
$ws.Cells.Item(71, 12).Value = 9683.7426
$ws.Cells.Item(71, 13).Value = -2486.584739999999
$ws.Cells.Item(71, 14).Value = -17795.7426

# CUL row 131
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(131, 8).Value = 2412.5232
$ws.Cells.Item(131, 10).Value = 2627.5195
$ws.Cells.Item(131, 12).Value = 7882.558499999999
$ws.Cells.Item(131, 14).Value = -17962.5585

# CUL row 132
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(132, 8).Value = 966.6667
$ws.Cells.Item(132, 10).Value = 1101.6
$ws.Cells.Item(132, 12).Value = 9914.4
$ws.Cells.Item(132, 14).Value = -14974.4

# GSM row 5
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(5, 8).Value = 1000
$ws.Cells.Item(5, 9).Value = 0
$ws.Cells.Item(5, 11).Value = 0
$ws.Cells.Item(5, 13).ClearContents()

# GSM row 117
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(117, 8).Value = 40000
$ws.Cells.Item(117, 10).Value = 40000
$ws.Cells.Item(117, 12).Value = 40000
$ws.Cells.Item(117, 14).Value = -46884

# GSM row 122
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(122, 8).Value = 1390251.4
$ws.Cells.Item(122, 9).Value = 11111111
$ws.Cells.Item(122, 10).Value = 1557.1428
$ws.Cells.Item(122, 11).Value = 33333333
$ws.Cells.Item(122, 12).Value = 4671.428400000001
$ws.Cells.Item(122, 13).Value = -33330883
$ws.Cells.Item(122, 14).Value = -9571.428400000001

# LTW row 2
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(2, 8).Value = 70002
$ws.Cells.Item(2, 10).Value = 70002
$ws.Cells.Item(2, 12).Value = 70002
$ws.Cells.Item(2, 14).Value = -70226

# LTW row 7
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(7, 8).Value = 3356.25
$ws.Cells.Item(7, 10).Value = 3360
$ws.Cells.Item(7, 12).Value = 3360
$ws.Cells.Item(7, 14).Value = -3584

# LTW row 69
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(69, 8).Value = 0
$ws.Cells.Item(69, 10).Value = 0
$ws.Cells.Item(69, 12).Value = 0
$ws.Cells.Item(69, 14).ClearContents()

# LTW row 72
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(72, 8).Value = 0
$ws.Cells.Item(72, 10).Value = 0
$ws.Cells.Item(72, 12).Value = 0
$ws.Cells.Item(72, 14).ClearContents()

# LTW row 118
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(118, 8).Value = 0
$ws.Cells.Item(118, 10).Value = 0
$ws.Cells.Item(118, 12).Value = 0
$ws.Cells.Item(118, 14).ClearContents()

# LTW row 126
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(126, 8).Value = 3356.25
$ws.Cells.Item(126, 10).Value = 3360
$ws.Cells.Item(126, 12).Value = 10080
$ws.Cells.Item(126, 14).Value = -15020

# LTW row 136
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(136, 8).Value = 3971.1562
$ws.Cells.Item(136, 9).Value = 1330.6818
$ws.Cells.Item(136, 10).Value = 9780.200000000001
$ws.Cells.Item(136, 11).Value = 3992.0454
$ws.Cells.Item(136, 12).Value = 29340.6
$ws.Cells.Item(136, 13).Value = -1442.0454
$ws.Cells.Item(136, 14).Value = -34440.60000000001

# WVR row 116
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(116, 8).Value = 20618
$ws.Cells.Item(116, 9).Value = 20618
$ws.Cells.Item(116, 10).Value = 0
$ws.Cells.Item(116, 11).Value = 20618
$ws.Cells.Item(116, 12).Value = 0
$ws.Cells.Item(116, 14).ClearContents()
$ws.Cells.Item(116, 13).Value = -16029

# WVR row 132
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(132, 8).Value = 10418206
$ws.Cells.Item(132, 9).Value = 12196257
$ws.Cells.Item(132, 10).Value = 3906.4285
$ws.Cells.Item(132, 11).Value = 36588771
$ws.Cells.Item(132, 12).Value = 11719.2855
$ws.Cells.Item(132, 13).Value = -36586241
$ws.Cells.Item(132, 14).Value = -16779.2855

# WVR row 136
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(136, 8).Value = 15198304
$ws.Cells.Item(136, 9).Value = 19667196
$ws.Cells.Item(136, 11).Value = 59001588
$ws.Cells.Item(136, 13).Value = -58999038
